$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 238.625
$ws.Range("I11").Value = 238.625
$ws.Range("K11").Value = 238.625
$ws.Range("M11").Value = -98.625
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H117").Value = 30885.143
$ws.Range("J117").Value = 30885.143
$ws.Range("L117").Value = 30885.143
$ws.Range("N117").Value = -40063.143
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H138").Value = 4351890
$ws.Range("J138").Value = 6901574.5
$ws.Range("L138").Value = 20704723.5
$ws.Range("N138").Value = -20715003.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 145.5
$ws.Range("I5").Value = 145.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 145.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -33.5
$ws.Range("N5").ClearContents()
$ws.Range("H45").Value = 3093.75
$ws.Range("I45").Value = 4087.5
$ws.Range("J45").Value = 2100
$ws.Range("K45").Value = 4087.5
$ws.Range("L45").Value = 2100
$ws.Range("M45").Value = -3710.5
$ws.Range("N45").Value = -2854
$ws.Range("H61").Value = 2813.0715
$ws.Range("I61").Value = 1632.5294
$ws.Range("J61").Value = 4637.5454
$ws.Range("K61").Value = 1632.5294
$ws.Range("L61").Value = 4637.5454
$ws.Range("M61").Value = -1420.5294
$ws.Range("N61").Value = -5061.5454
$ws.Range("H111").Value = 31111.223
$ws.Range("J111").Value = 31111.223
$ws.Range("L111").Value = 31111.223
$ws.Range("N111").Value = -39291.223
$ws.Range("H113").Value = 30000
$ws.Range("J113").Value = 30000
$ws.Range("L113").Value = 30000
$ws.Range("N113").Value = -38678
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H136").Value = 2813.0715
$ws.Range("I136").Value = 1632.5294
$ws.Range("J136").Value = 4637.5454
$ws.Range("K136").Value = 4897.5882
$ws.Range("L136").Value = 13912.6362
$ws.Range("M136").Value = -2347.5882
$ws.Range("N136").Value = -19012.6362
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 145.5
$ws.Range("I4").Value = 145.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 145.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -30.5
$ws.Range("N4").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 105
$ws.Range("I7").Value = 56.11111
$ws.Range("J7").Value = 145
$ws.Range("K7").Value = 56.11111
$ws.Range("L7").Value = 145
$ws.Range("M7").Value = 56.88889
$ws.Range("N7").Value = -371
$ws.Range("H22").Value = 533.2143
$ws.Range("I22").Value = 424.1111
$ws.Range("J22").Value = 729.6
$ws.Range("K22").Value = 424.1111
$ws.Range("L22").Value = 729.6
$ws.Range("M22").Value = -74.11110000000002
$ws.Range("N22").Value = -1429.6
$ws.Range("H31").Value = 8475654
$ws.Range("I31").Value = 776.2941
$ws.Range("J31").Value = 11905962
$ws.Range("K31").Value = 776.2941
$ws.Range("L31").Value = 11905962
$ws.Range("M31").Value = -481.2941
$ws.Range("N31").Value = -11906552
$ws.Range("H34").Value = 8475654
$ws.Range("I34").Value = 776.2941
$ws.Range("J34").Value = 11905962
$ws.Range("K34").Value = 776.2941
$ws.Range("L34").Value = 11905962
$ws.Range("M34").Value = -574.2941
$ws.Range("N34").Value = -11906366
$ws.Range("H58").Value = 1863.6666
$ws.Range("I58").Value = 1846.625
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1846.625
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -1643.625
$ws.Range("N58").Value = -2406
$ws.Range("H105").Value = 1413.4286
$ws.Range("I105").Value = 886.25
$ws.Range("J105").Value = 2116.3333
$ws.Range("K105").Value = 886.25
$ws.Range("L105").Value = 2116.3333
$ws.Range("M105").Value = 860.75
$ws.Range("N105").Value = -5610.3333
$ws.Range("H132").Value = 3726.68
$ws.Range("I132").Value = 3530.8333
$ws.Range("J132").Value = 4230.2856
$ws.Range("K132").Value = 10592.4999
$ws.Range("L132").Value = 12690.8568
$ws.Range("M132").Value = -8062.499899999999
$ws.Range("N132").Value = -17750.8568
$ws.Range("H134").Value = 2334.2415
$ws.Range("I134").Value = 2201.1428
$ws.Range("J134").Value = 2683.625
$ws.Range("K134").Value = 6603.428400000001
$ws.Range("L134").Value = 8050.875
$ws.Range("M134").Value = -4068.428400000001
$ws.Range("N134").Value = -13120.875
$ws.Range("H136").Value = 1863.6666
$ws.Range("I136").Value = 1846.625
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5539.875
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2989.875
$ws.Range("N136").Value = -11100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2958.1052
$ws.Range("I51").Value = 499.75
$ws.Range("J51").Value = 3613.6667
$ws.Range("K51").Value = 1499.25
$ws.Range("L51").Value = 10841.0001
$ws.Range("M51").Value = -1039.25
$ws.Range("N51").Value = -11761.0001
$ws.Range("H68").Value = 1885.75
$ws.Range("I68").Value = 1061
$ws.Range("J68").Value = 2276.4211
$ws.Range("K68").Value = 3183
$ws.Range("L68").Value = 6829.263300000001
$ws.Range("M68").Value = -2372
$ws.Range("N68").Value = -8451.263300000001
$ws.Range("H71").Value = 1885.75
$ws.Range("I71").Value = 1061
$ws.Range("J71").Value = 2276.4211
$ws.Range("K71").Value = 9549
$ws.Range("L71").Value = 20487.7899
$ws.Range("M71").Value = -5493
$ws.Range("N71").Value = -28599.7899
$ws.Range("H131").Value = 880.22
$ws.Range("I131").Value = 617.75
$ws.Range("J131").Value = 891.15625
$ws.Range("K131").Value = 1853.25
$ws.Range("L131").Value = 2673.46875
$ws.Range("M131").Value = 3186.75
$ws.Range("N131").Value = -12753.46875
$ws.Range("H133").Value = 8750
$ws.Range("I133").Value = 8750
$ws.Range("K133").Value = 26250
$ws.Range("M133").Value = -21190
$ws.Range("H137").Value = 18293630
$ws.Range("I137").Value = 2280.9092
$ws.Range("J137").Value = 26032278
$ws.Range("K137").Value = 6842.7276
$ws.Range("L137").Value = 78096834
$ws.Range("M137").Value = -1742.7276
$ws.Range("N137").Value = -78107034
$ws.Range("H139").Value = 1847.5
$ws.Range("I139").Value = 1226.9231
$ws.Range("J139").Value = 3000
$ws.Range("K139").Value = 3680.7693
$ws.Range("L139").Value = 9000
$ws.Range("M139").Value = 1459.2307
$ws.Range("N139").Value = -19280
$ws.Range("H141").Value = 1140
$ws.Range("I141").Value = 1125
$ws.Range("K141").Value = 3375
$ws.Range("M141").Value = 1805
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 746.4167
$ws.Range("I22").Value = 445.5
$ws.Range("K22").Value = 445.5
$ws.Range("M22").Value = -150.5
$ws.Range("H27").Value = 746.4167
$ws.Range("I27").Value = 445.5
$ws.Range("K27").Value = 445.5
$ws.Range("M27").Value = -338.5
$ws.Range("H93").Value = 72961.60000000001
$ws.Range("I93").Value = 3326
$ws.Range("J93").Value = 351504
$ws.Range("K93").Value = 3326
$ws.Range("L93").Value = 351504
$ws.Range("M93").Value = -2078
$ws.Range("N93").Value = -354000
$ws.Range("H132").Value = 3159.2678
$ws.Range("I132").Value = 2674.743
$ws.Range("J132").Value = 3966.8096
$ws.Range("K132").Value = 8024.228999999999
$ws.Range("L132").Value = 11900.4288
$ws.Range("M132").Value = -5494.228999999999
$ws.Range("N132").Value = -16960.4288
$ws.Range("H136").Value = 1722.238
$ws.Range("I136").Value = 1458.8
$ws.Range("J136").Value = 2380.8333
$ws.Range("K136").Value = 4376.4
$ws.Range("L136").Value = 7142.499899999999
$ws.Range("M136").Value = -1826.4
$ws.Range("N136").Value = -12242.4999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2553.0513
$ws.Range("I136").Value = 2855.04
$ws.Range("J136").Value = 2013.7858
$ws.Range("K136").Value = 8565.119999999999
$ws.Range("L136").Value = 6041.357400000001
$ws.Range("M136").Value = -6015.119999999999
$ws.Range("N136").Value = -11141.3574
